$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new boolean value TRUE in cell E3
$ws.Range("E3").Value = $true

# Update the active selection to E3, as seen after the edit
$ws.Range("E3").Select()
